# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 17:35"

# --- Simple data refresh for Gibraltar (row 159): B, C, E columns ---
$ws.Range("B159").Value = 149
$ws.Range("C159").Value = 2
$ws.Range("E159").Value = 4

# --- Simple data refresh for Libia (row 176): B, C, E columns ---
$ws.Range("B176").Value = 69
$ws.Range("C176").Value = 1
$ws.Range("E176").Value = 31

# --- Swap Santa Lucia (row 196) and Belice (row 197), along with their D/H data ---
$ws.Range("A196").Value = "Belice"
$ws.Range("D196").Value = 16
$ws.Range("H196").Value = 2

$ws.Range("A197").Value = "Santa Lucia"
$ws.Range("D197").Value = 18
$ws.Range("H197").Value = 0

# --- Swap Montserrat (row 209) and Groenlandia (row 210), along with their D/H data ---
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1
